# feat: add 2022-Q4 data
#
# Plan:
#  1. Duplicate the existing "2022-Q3" sheet (placing the copy *after* it).
#     The original keeps its sheetId/rId; we will turn the original into the
#     new "2022-Q4" sheet (fresh data) and leave the copy as the untouched
#     "2022-Q3" sheet (old data, unchanged) - this reproduces the exact
#     sheetId/r:id numbering from the target workbook.
#  2. Rename the two sheets accordingly.
#  3. Overwrite the (renamed) "2022-Q4" sheet's fund table with the new
#     2022-Q4 figures.
#  4. Insert the 2022-Q4 summary row on the "总计" sheet, pushing the old
#     2022-Q3 summary row down to row 3.

$wb = $excel.ActiveWorkbook

function Set-TextValue($rng, $txt) {
    $rng.NumberFormat = "@"
    $rng.Value = $txt
    $rng.Style = "Normal"
}

# ---------------------------------------------------------------------
# Step 1+2: duplicate "2022-Q3", rename copies
# ---------------------------------------------------------------------
$wsOldQ3 = $wb.Worksheets.Item("2022-Q3")
$wsOldQ3.Copy($null, $wsOldQ3)

$wsQ4 = $wb.Worksheets.Item(2)
$wsQ3 = $wb.Worksheets.Item(3)
$wsQ4.Name = "2022-Q4"
$wsQ3.Name = "2022-Q3"

# ---------------------------------------------------------------------
# Step 3: repopulate the "2022-Q4" sheet with the new fund table
# ---------------------------------------------------------------------
$fundRows = @(
    @(0, "000800", "华商未来主题混合", "4.12", "74.31", "3.17", "0.1306", 8),
    @(1, "010656", "华商均衡30混合", "3.86", "33.52", "2.23", "0.0861", 7),
    @(2, "460009", "华泰柏瑞量化先行混合A", "4.35", "94.47", "0.98", "0.0426", 4),
    @(3, "002289", "华商改革创新股票A", "1.12", "85.85", "3.23", "0.0362", 9),
    @(4, "005055", "华泰柏瑞量化阿尔法灵活配置混合A", "2.10", "93.30", "1.04", "0.0218", 9),
    @(5, "010403", "华商景气优选混合", "0.58", "77.20", "3.69", "0.0214", 8),
    @(6, "016052", "华商改革创新股票C", "0.48", "85.85", "3.23", "0.0155", 9),
    @(7, "010246", "华泰柏瑞量化先行混合C", "0.77", "94.47", "0.98", "0.0075", 4),
    @(8, "006532", "华泰柏瑞量化阿尔法灵活配置混合C", "0.02", "93.30", "1.04", "0.0002", 9)
)

$fundRowNum = 2
foreach ($fundRow in $fundRows) {
    $wsQ4.Range("A$fundRowNum").Value = $fundRow[0]
    Set-TextValue $wsQ4.Range("B$fundRowNum") $fundRow[1]
    Set-TextValue $wsQ4.Range("C$fundRowNum") $fundRow[2]
    Set-TextValue $wsQ4.Range("D$fundRowNum") $fundRow[3]
    Set-TextValue $wsQ4.Range("E$fundRowNum") $fundRow[4]
    Set-TextValue $wsQ4.Range("F$fundRowNum") $fundRow[5]
    Set-TextValue $wsQ4.Range("G$fundRowNum") $fundRow[6]
    $wsQ4.Range("H$fundRowNum").Value = $fundRow[7]
    $fundRowNum = $fundRowNum + 1
}

# New table only has 9 data rows (2..10); the old table had 7 (2..8), so
# nothing needs clearing - every previously used row is overwritten and the
# two new rows (9, 10) are brand new.

# Match header/column-A styling ("s=2") to the rest of the workbook by
# copying the already-correctly-styled format from the "总计" sheet.
$wsTotal = $wb.Worksheets.Item("总计")
$wsTotal.Range("B1").Copy()
$wsQ4.Range("B1:H1").PasteSpecial(-4122)
$wsTotal.Range("A2").Copy()
$wsQ4.Range("A2:A10").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Step 4: update the "总计" sheet - insert the 2022-Q4 summary row
# ---------------------------------------------------------------------
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q3"
$wsTotal.Range("C3").Value = 7
$wsTotal.Range("D3").Value = 0.29
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)

$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 9
$wsTotal.Range("D2").Value = 0.36

$wsTotal.Range("A1").Select()
